$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the report date: set J2 to the user-entered date (no time component),
# and change its number format to a date-only format instead of date+time.
$ws.Range("J2").Value = 44393
$ws.Range("J2").NumberFormat = "m/d/yy"

# Fix divide-by-zero / unassigned variable issues: populate Quantity/Total/Average
# values that were previously left at 0 / 65535 sentinel values.

# Store Summary
$ws.Range("B6").Value = 20
$ws.Range("C6").Value = 957.58
$ws.Range("D6").Value = 47.879000000000005

$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 48829.65
$ws.Range("D7").Value = 12207.4125

$ws.Range("B8").Value = 20
$ws.Range("C8").Value = 900.81
$ws.Range("D8").Value = 45.040499999999994

$ws.Range("C9").Value = 56.77

$ws.Range("B10").Value = 20
$ws.Range("C10").Value = 957.57999999999993
$ws.Range("D10").Value = 47.878999999999998

# Payment Types
$ws.Range("B13").Value = 12
$ws.Range("C13").Value = 585.45000000000005
$ws.Range("D13").Value = 585.45000000000005

$ws.Range("B14").Value = 12
$ws.Range("C14").Value = 49201.78
$ws.Range("D14").Value = 49201.78

# Tax Categories
$ws.Range("C21").Value = 900.81
